$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.548.43"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "2.426.14"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "564.45"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").Value = "166.39"
$ws.Range("E6").Value = "  +5.57%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("E9").Value = "  +8.16%  "
$ws.Range("D10").Value = "2.424.18"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").Value = "69.304.77"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("E15").Value = "  +5.13%  "
$ws.Range("D16").Value = "2.871.14"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "23.99"
$ws.Range("E17").Value = "  +5.17%  "
$ws.Range("D18").Value = "2.430.45"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "10.80"
$ws.Range("E19").Value = "  +4.61%  "
$ws.Range("D20").Value = "342.47"
$ws.Range("E20").Value = "  +4.37%  "
$ws.Range("E21").Value = "  +5.09%  "
$ws.Range("D22").Value = "3.88"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("E23").Value = "  +6.97%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "66.19"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("E26").Value = "  +6.21%  "
$ws.Range("D27").Value = "8.50"
$ws.Range("E27").Value = "  +6.59%  "
$ws.Range("D28").Value = "2.554.83"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "0.0₃0849"
$ws.Range("E30").Value = "  +6.53%  "
$ws.Range("D31").Value = "7.40"
$ws.Range("E31").Value = "  +6.25%  "
$ws.Range("D32").Value = "1.24"
$ws.Range("E32").Value = "  +11.66%  "
$ws.Range("D33").Value = "454.00"
$ws.Range("E33").Value = "  +9.72%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("D36").Value = "158.89"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("E38").Value = "  +6.56%  "
$ws.Range("D40").Value = "18.22"
$ws.Range("E40").Value = "  +2.71%  "
$ws.Range("D41").Value = "0.303"
$ws.Range("E41").Value = "  +3.79%  "
$ws.Range("E42").Value = "  +5.38%  "
$ws.Range("D43").Value = "4.40"
$ws.Range("E43").Value = "  +4.80%  "
$ws.Range("D44").Value = "37.84"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("E45").Value = "  +2.60%  "
$ws.Range("E46").Value = "  +6.22%  "
$ws.Range("D47").Value = "135.03"
$ws.Range("E47").Value = "  +4.73%  "
$ws.Range("D48").Value = "3.40"
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("E50").Value = "  +3.16%  "
$ws.Range("D51").Value = "0.0935"
$ws.Range("E51").Value = "  +2.70%  "
